$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark (currently sitting at the end
#    of the "You can also embed plots, for example:" paragraph).
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. Re-create "_GoBack" at the very start of the document (start of
#    the Title paragraph, before the "Bayesian Calibration" run) so
#    it becomes bookmark id 0, pushing "r-markdown" to id 1 and
#    "including-plots" to id 2 -- matching the target diff.
#
#    A bookmark collapsed exactly at absolute document offset 0 tends
#    to snap to cover the whole first paragraph, so we insert a
#    throw-away character at offset 0, anchor the bookmark right
#    after it (offset 1, a perfectly ordinary collapsed position),
#    and then delete the throw-away character again. The bookmark
#    stays put and ends up collapsed at true offset 0.
# ------------------------------------------------------------------
$tmp = $d.Range(0, 0)
$tmp.InsertBefore("Z")

$atStart = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $atStart)

$d.Range(0, 1).Delete()

# ------------------------------------------------------------------
# 3. Fix the "Table Caption" style so its run properties no longer
#    inherit italics from the base "Caption" style.
# ------------------------------------------------------------------
$tableCaption = $d.Styles("TableCaption")
$tableCaption.Font.Italic = 0

Write-Output "done"
